$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, matching the style/format of the existing header row (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill the new "Save" column with 0 for each data row
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
